$d = $word.ActiveDocument

# Locate the second "Ruta de tu casa" paragraph (the one right before the
# trailing empty paragraph) and insert the new "LANG" section after it,
# mirroring the existing "HOME" section.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Ruta de tu casa") {
        $anchor = $p
    }
}

$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$headingPara = $anchor.Next()
$headingPara.Style = "Heading 1"
$headingPara.Range.Text = "LANG"

$hr = $headingPara.Range
$hr.Collapse(0)
$hr.InsertParagraphAfter()

$descPara = $headingPara.Next()
$descPara.Style = "Normal"
$descPara.Range.Text = "Mostrar tu lenguaje y sistema de codificación."
